$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 50.0
$ws.Range("D3").Value = 98.0
$ws.Range("D4").Value = 121
$ws.Range("D5").Value = 59
$ws.Range("D6").Value = 11.0

$ws.Range("E6").Select()
